$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (35 cell updates) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 2340  # H17: 1399.375 -> 2340
$ws.Cells.Item(17, 10).Value = 3333.3333  # J17: 1582.5 -> 3333.3333
$ws.Cells.Item(17, 12).Value = 9999.999899999999  # L17: 4747.5 -> 9999.999899999999
$ws.Cells.Item(17, 14).Value = -10335.9999  # N17: -5083.5 -> -10335.9999
$ws.Cells.Item(18, 8).Value = 1018.4  # H18: 1137.3846 -> 1018.4
$ws.Cells.Item(18, 9).Value = 876.8570999999999  # I18: 982.1667 -> 876.8570999999999
$ws.Cells.Item(18, 11).Value = 876.8570999999999  # K18: 982.1667 -> 876.8570999999999
$ws.Cells.Item(18, 13).Value = -592.8570999999999  # M18: -698.1667 -> -592.8570999999999
$ws.Cells.Item(112, 8).Value = 930.7778  # H112: 929.6429000000001 -> 930.7778
$ws.Cells.Item(112, 9).Value = 683.6667  # I112: 714.4286 -> 683.6667
$ws.Cells.Item(112, 11).Value = 2051.0001  # K112: 2143.2858 -> 2051.0001
$ws.Cells.Item(112, 13).Value = -943.0001000000002  # M112: -1035.2858 -> -943.0001000000002
$ws.Cells.Item(132, 8).Value = 1851.5385  # H132: 1877.6296 -> 1851.5385
$ws.Cells.Item(132, 9).Value = 1885.6  # I132: 2122.7273 -> 1885.6
$ws.Cells.Item(132, 10).Value = 1000  # J132: 799.2 -> 1000
$ws.Cells.Item(132, 11).Value = 5656.799999999999  # K132: 6368.1819 -> 5656.799999999999
$ws.Cells.Item(132, 12).Value = 3000  # L132: 2397.6 -> 3000
$ws.Cells.Item(132, 13).Value = -3126.799999999999  # M132: -3838.1819 -> -3126.799999999999
$ws.Cells.Item(132, 14).Value = -8060  # N132: -7457.6 -> -8060
$ws.Cells.Item(133, 8).Value = 79329  # H133: 78526.61 -> 79329
$ws.Cells.Item(133, 10).Value = 79329  # J133: 78526.61 -> 79329
$ws.Cells.Item(133, 12).Value = 79329  # L133: 78526.61 -> 79329
$ws.Cells.Item(133, 14).Value = -89449  # N133: -88646.61 -> -89449
$ws.Cells.Item(134, 8).Value = 94216.664  # H134: 89042 -> 94216.664
$ws.Cells.Item(134, 10).Value = 94216.664  # J134: 89042 -> 94216.664
$ws.Cells.Item(134, 12).Value = 94216.664  # L134: 89042 -> 94216.664
$ws.Cells.Item(134, 14).Value = -104356.664  # N134: -99182 -> -104356.664
$ws.Cells.Item(137, 8).Value = 280947.84  # H137: 275663.72 -> 280947.84
$ws.Cells.Item(137, 9).Value = 1982.9688  # I137: 1949.8182 -> 1982.9688
$ws.Cells.Item(137, 11).Value = 5948.9064  # K137: 5849.4546 -> 5948.9064
$ws.Cells.Item(137, 13).Value = -3398.9064  # M137: -3299.4546 -> -3398.9064
$ws.Cells.Item(138, 8).Value = 1463.0212  # H138: 1439.6522 -> 1463.0212
$ws.Cells.Item(138, 10).Value = 2623.2727  # J138: 2631.8 -> 2623.2727
$ws.Cells.Item(138, 12).Value = 7869.8181  # L138: 7895.400000000001 -> 7869.8181
$ws.Cells.Item(138, 14).Value = -18149.8181  # N138: -18175.4 -> -18149.8181

# ---- Sheet: ARM (15 cell updates) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 545.4  # H2: 661.875 -> 545.4
$ws.Cells.Item(2, 9).Value = 439.33334  # I2: 542.1429000000001 -> 439.33334
$ws.Cells.Item(2, 11).Value = 439.33334  # K2: 542.1429000000001 -> 439.33334
$ws.Cells.Item(2, 13).Value = -326.33334  # M2: -429.1429000000001 -> -326.33334
$ws.Cells.Item(116, 8).Value = 545.4  # H116: 661.875 -> 545.4
$ws.Cells.Item(116, 9).Value = 439.33334  # I116: 542.1429000000001 -> 439.33334
$ws.Cells.Item(116, 11).Value = 439.33334  # K116: 542.1429000000001 -> 439.33334
$ws.Cells.Item(116, 13).Value = 1854.66666  # M116: 1751.8571 -> 1854.66666
$ws.Cells.Item(132, 8).Value = 1585.8823  # H132: 1294.575 -> 1585.8823
$ws.Cells.Item(132, 9).Value = 1089.6666  # I132: 914.8182 -> 1089.6666
$ws.Cells.Item(132, 10).Value = 3499.8572  # J132: 3084.8572 -> 3499.8572
$ws.Cells.Item(132, 11).Value = 3268.9998  # K132: 2744.4546 -> 3268.9998
$ws.Cells.Item(132, 12).Value = 10499.5716  # L132: 9254.571599999999 -> 10499.5716
$ws.Cells.Item(132, 13).Value = -738.9998000000001  # M132: -214.4546 -> -738.9998000000001
$ws.Cells.Item(132, 14).Value = -15559.5716  # N132: -14314.5716 -> -15559.5716

# ---- Sheet: BSM (12 cell updates) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 545.4  # H3: 661.875 -> 545.4
$ws.Cells.Item(3, 9).Value = 439.33334  # I3: 542.1429000000001 -> 439.33334
$ws.Cells.Item(3, 11).Value = 439.33334  # K3: 542.1429000000001 -> 439.33334
$ws.Cells.Item(3, 13).Value = -325.33334  # M3: -428.1429000000001 -> -325.33334
$ws.Cells.Item(134, 8).Value = 1110.125  # H134: 1117.8572 -> 1110.125
$ws.Cells.Item(134, 9).Value = 1110.125  # I134: 1117.8572 -> 1110.125
$ws.Cells.Item(134, 11).Value = 3330.375  # K134: 3353.5716 -> 3330.375
$ws.Cells.Item(134, 13).Value = -795.375  # M134: -818.5715999999998 -> -795.375
$ws.Cells.Item(135, 8).Value = 45584.633  # H135: 42302.36 -> 45584.633
$ws.Cells.Item(135, 10).Value = 45584.633  # J135: 42302.36 -> 45584.633
$ws.Cells.Item(135, 12).Value = 45584.633  # L135: 42302.36 -> 45584.633
$ws.Cells.Item(135, 14).Value = -55724.633  # N135: -52442.36 -> -55724.633

# ---- Sheet: CRP (54 cell updates) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2028.6857  # H31: 2097.121 -> 2028.6857
$ws.Cells.Item(31, 9).Value = 1693.7587  # I31: 1725.7142 -> 1693.7587
$ws.Cells.Item(31, 10).Value = 3647.5  # J31: 4177 -> 3647.5
$ws.Cells.Item(31, 11).Value = 1693.7587  # K31: 1725.7142 -> 1693.7587
$ws.Cells.Item(31, 12).Value = 3647.5  # L31: 4177 -> 3647.5
$ws.Cells.Item(31, 13).Value = -1398.7587  # M31: -1430.7142 -> -1398.7587
$ws.Cells.Item(31, 14).Value = -4237.5  # N31: -4767 -> -4237.5
$ws.Cells.Item(34, 8).Value = 2028.6857  # H34: 2097.121 -> 2028.6857
$ws.Cells.Item(34, 9).Value = 1693.7587  # I34: 1725.7142 -> 1693.7587
$ws.Cells.Item(34, 10).Value = 3647.5  # J34: 4177 -> 3647.5
$ws.Cells.Item(34, 11).Value = 1693.7587  # K34: 1725.7142 -> 1693.7587
$ws.Cells.Item(34, 12).Value = 3647.5  # L34: 4177 -> 3647.5
$ws.Cells.Item(34, 13).Value = -1491.7587  # M34: -1523.7142 -> -1491.7587
$ws.Cells.Item(34, 14).Value = -4051.5  # N34: -4581 -> -4051.5
$ws.Cells.Item(58, 8).Value = 1375.5272  # H58: 1432.1765 -> 1375.5272
$ws.Cells.Item(58, 9).Value = 1322.1714  # I58: 1356.8823 -> 1322.1714
$ws.Cells.Item(58, 10).Value = 1468.9  # J58: 1582.7646 -> 1468.9
$ws.Cells.Item(58, 11).Value = 1322.1714  # K58: 1356.8823 -> 1322.1714
$ws.Cells.Item(58, 12).Value = 1468.9  # L58: 1582.7646 -> 1468.9
$ws.Cells.Item(58, 13).Value = -1119.1714  # M58: -1153.8823 -> -1119.1714
$ws.Cells.Item(58, 14).Value = -1874.9  # N58: -1988.7646 -> -1874.9
$ws.Cells.Item(86, 8).Value = 4397.1724  # H86: 4636.2856 -> 4397.1724
$ws.Cells.Item(86, 9).Value = 3409.889  # I86: 3745.647 -> 3409.889
$ws.Cells.Item(86, 11).Value = 3409.889  # K86: 3745.647 -> 3409.889
$ws.Cells.Item(86, 13).Value = -2286.889  # M86: -2622.647 -> -2286.889
$ws.Cells.Item(89, 8).Value = 4397.1724  # H89: 4636.2856 -> 4397.1724
$ws.Cells.Item(89, 9).Value = 3409.889  # I89: 3745.647 -> 3409.889
$ws.Cells.Item(89, 11).Value = 17049.445  # K89: 18728.235 -> 17049.445
$ws.Cells.Item(89, 13).Value = -11433.445  # M89: -13112.235 -> -11433.445
$ws.Cells.Item(122, 8).Value = 2286.25  # H122: 2437.1155 -> 2286.25
$ws.Cells.Item(122, 9).Value = 2067.375  # I122: 2205.7856 -> 2067.375
$ws.Cells.Item(122, 10).Value = 2578.0833  # J122: 2707 -> 2578.0833
$ws.Cells.Item(122, 11).Value = 6202.125  # K122: 6617.3568 -> 6202.125
$ws.Cells.Item(122, 12).Value = 7734.249899999999  # L122: 8121 -> 7734.249899999999
$ws.Cells.Item(122, 13).Value = -3752.125  # M122: -4167.3568 -> -3752.125
$ws.Cells.Item(122, 14).Value = -12634.2499  # N122: -13021 -> -12634.2499
$ws.Cells.Item(125, 8).Value = 59500  # H125: 0 -> 59500
$ws.Cells.Item(125, 10).Value = 59500  # J125: 0 -> 59500
$ws.Cells.Item(125, 12).Value = 59500  # L125: 0 -> 59500
$ws.Cells.Item(125, 14).Value = -64420  # N125: None -> -64420
$ws.Cells.Item(134, 8).Value = 64479.875  # H134: 79055.69500000001 -> 64479.875
$ws.Cells.Item(134, 9).Value = 1362.5714  # I134: 1516.8 -> 1362.5714
$ws.Cells.Item(134, 10).Value = 113571.11  # J134: 127517.5 -> 113571.11
$ws.Cells.Item(134, 11).Value = 4087.7142  # K134: 4550.4 -> 4087.7142
$ws.Cells.Item(134, 12).Value = 340713.33  # L134: 382552.5 -> 340713.33
$ws.Cells.Item(134, 13).Value = -1552.7142  # M134: -2015.4 -> -1552.7142
$ws.Cells.Item(134, 14).Value = -345783.33  # N134: -387622.5 -> -345783.33
$ws.Cells.Item(136, 8).Value = 1375.5272  # H136: 1432.1765 -> 1375.5272
$ws.Cells.Item(136, 9).Value = 1322.1714  # I136: 1356.8823 -> 1322.1714
$ws.Cells.Item(136, 10).Value = 1468.9  # J136: 1582.7646 -> 1468.9
$ws.Cells.Item(136, 11).Value = 3966.5142  # K136: 4070.6469 -> 3966.5142
$ws.Cells.Item(136, 12).Value = 4406.700000000001  # L136: 4748.293799999999 -> 4406.700000000001
$ws.Cells.Item(136, 13).Value = -1416.5142  # M136: -1520.6469 -> -1416.5142
$ws.Cells.Item(136, 14).Value = -9506.700000000001  # N136: -9848.293799999999 -> -9506.700000000001

# ---- Sheet: CUL (20 cell updates) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(48, 8).Value = 5986  # H48: 5818.3335 -> 5986
$ws.Cells.Item(48, 10).Value = 5986  # J48: 5818.3335 -> 5986
$ws.Cells.Item(48, 12).Value = 17958  # L48: 17455.0005 -> 17958
$ws.Cells.Item(48, 14).Value = -18458  # N48: -17955.0005 -> -18458
$ws.Cells.Item(56, 8).Value = 6993.5  # H56: 6994 -> 6993.5
$ws.Cells.Item(56, 9).Value = 6993.5  # I56: 6994 -> 6993.5
$ws.Cells.Item(56, 11).Value = 6993.5  # K56: 6994 -> 6993.5
$ws.Cells.Item(56, 13).Value = -6463.5  # M56: -6464 -> -6463.5
$ws.Cells.Item(107, 8).Value = 395.2857  # H107: 396.25 -> 395.2857
$ws.Cells.Item(107, 10).Value = 395.2857  # J107: 396.25 -> 395.2857
$ws.Cells.Item(107, 12).Value = 1185.8571  # L107: 1188.75 -> 1185.8571
$ws.Cells.Item(107, 14).Value = -5025.8571  # N107: -5028.75 -> -5025.8571
$ws.Cells.Item(122, 8).Value = 497.55554  # H122: 522.25 -> 497.55554
$ws.Cells.Item(122, 10).Value = 482.7143  # J122: 513.1667 -> 482.7143
$ws.Cells.Item(122, 12).Value = 4344.428699999999  # L122: 4618.5003 -> 4344.428699999999
$ws.Cells.Item(122, 14).Value = -9244.4287  # N122: -9518.5003 -> -9244.4287
$ws.Cells.Item(124, 8).Value = 0  # H124: 1500 -> 0
$ws.Cells.Item(124, 9).Value = 0  # I124: 1500 -> 0
$ws.Cells.Item(124, 11).Value = 0  # K124: 4500 -> 0
$ws.Cells.Item(124, 13).ClearContents()  # M124 was 410, now blank

# ---- Sheet: GSM (34 cell updates) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(54, 8).Value = 2898.75  # H54: 3598.75 -> 2898.75
$ws.Cells.Item(54, 10).Value = 2898.75  # J54: 3598.75 -> 2898.75
$ws.Cells.Item(54, 12).Value = 2898.75  # L54: 3598.75 -> 2898.75
$ws.Cells.Item(54, 14).Value = -3678.75  # N54: -4378.75 -> -3678.75
$ws.Cells.Item(70, 8).Value = 7172.875  # H70: 7396.077 -> 7172.875
$ws.Cells.Item(70, 9).Value = 7208.778  # I70: 7467.2856 -> 7208.778
$ws.Cells.Item(70, 10).Value = 7126.7144  # J70: 7313 -> 7126.7144
$ws.Cells.Item(70, 11).Value = 7208.778  # K70: 7467.2856 -> 7208.778
$ws.Cells.Item(70, 12).Value = 7126.7144  # L70: 7313 -> 7126.7144
$ws.Cells.Item(70, 13).Value = -6938.778  # M70: -7197.2856 -> -6938.778
$ws.Cells.Item(70, 14).Value = -7666.7144  # N70: -7853 -> -7666.7144
$ws.Cells.Item(73, 8).Value = 7172.875  # H73: 7396.077 -> 7172.875
$ws.Cells.Item(73, 9).Value = 7208.778  # I73: 7467.2856 -> 7208.778
$ws.Cells.Item(73, 10).Value = 7126.7144  # J73: 7313 -> 7126.7144
$ws.Cells.Item(73, 11).Value = 7208.778  # K73: 7467.2856 -> 7208.778
$ws.Cells.Item(73, 12).Value = 7126.7144  # L73: 7313 -> 7126.7144
$ws.Cells.Item(73, 13).Value = -6272.778  # M73: -6531.2856 -> -6272.778
$ws.Cells.Item(73, 14).Value = -8998.714400000001  # N73: -9185 -> -8998.714400000001
$ws.Cells.Item(92, 8).Value = 9116.5  # H92: 9364.625 -> 9116.5
$ws.Cells.Item(92, 10).Value = 9116.5  # J92: 9364.625 -> 9116.5
$ws.Cells.Item(92, 12).Value = 9116.5  # L92: 9364.625 -> 9116.5
$ws.Cells.Item(92, 14).Value = -12860.5  # N92: -13108.625 -> -12860.5
$ws.Cells.Item(102, 8).Value = 2155.2778  # H102: 2223.4119 -> 2155.2778
$ws.Cells.Item(102, 9).Value = 1992.7858  # I102: 2069.3845 -> 1992.7858
$ws.Cells.Item(102, 11).Value = 1992.7858  # K102: 2069.3845 -> 1992.7858
$ws.Cells.Item(102, 13).Value = -370.7858000000001  # M102: -447.3845000000001 -> -370.7858000000001
$ws.Cells.Item(131, 8).Value = 98664.664  # H131: 98797.8 -> 98664.664
$ws.Cells.Item(131, 10).Value = 98664.664  # J131: 98797.8 -> 98664.664
$ws.Cells.Item(131, 12).Value = 98664.664  # L131: 98797.8 -> 98664.664
$ws.Cells.Item(131, 14).Value = -108744.664  # N131: -108877.8 -> -108744.664
$ws.Cells.Item(132, 8).Value = 4568.5356  # H132: 4207.4517 -> 4568.5356
$ws.Cells.Item(132, 9).Value = 3089.5454  # I132: 2819.28 -> 3089.5454
$ws.Cells.Item(132, 11).Value = 9268.636200000001  # K132: 8457.84 -> 9268.636200000001
$ws.Cells.Item(132, 13).Value = -6738.636200000001  # M132: -5927.84 -> -6738.636200000001

# ---- Sheet: LTW (44 cell updates) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(5, 8).Value = 50000  # H5: 0 -> 50000
$ws.Cells.Item(5, 10).Value = 50000  # J5: 0 -> 50000
$ws.Cells.Item(5, 12).Value = 50000  # L5: 0 -> 50000
$ws.Cells.Item(5, 14).Value = -50226  # N5: None -> -50226
$ws.Cells.Item(22, 8).Value = 4600.654  # H22: 4289.607 -> 4600.654
$ws.Cells.Item(22, 9).Value = 721.53845  # I22: 676.7857 -> 721.53845
$ws.Cells.Item(22, 10).Value = 8479.77  # J22: 7902.4287 -> 8479.77
$ws.Cells.Item(22, 11).Value = 721.53845  # K22: 676.7857 -> 721.53845
$ws.Cells.Item(22, 12).Value = 8479.77  # L22: 7902.4287 -> 8479.77
$ws.Cells.Item(22, 13).Value = -426.53845  # M22: -381.7857 -> -426.53845
$ws.Cells.Item(22, 14).Value = -9069.77  # N22: -8492.4287 -> -9069.77
$ws.Cells.Item(27, 8).Value = 4600.654  # H27: 4289.607 -> 4600.654
$ws.Cells.Item(27, 9).Value = 721.53845  # I27: 676.7857 -> 721.53845
$ws.Cells.Item(27, 10).Value = 8479.77  # J27: 7902.4287 -> 8479.77
$ws.Cells.Item(27, 11).Value = 721.53845  # K27: 676.7857 -> 721.53845
$ws.Cells.Item(27, 12).Value = 8479.77  # L27: 7902.4287 -> 8479.77
$ws.Cells.Item(27, 13).Value = -614.53845  # M27: -569.7857 -> -614.53845
$ws.Cells.Item(27, 14).Value = -8693.77  # N27: -8116.4287 -> -8693.77
$ws.Cells.Item(46, 8).Value = 8657.177  # H46: 9017 -> 8657.177
$ws.Cells.Item(46, 9).Value = 12149.6  # I46: 13388.444 -> 12149.6
$ws.Cells.Item(46, 10).Value = 3668  # J46: 3396.5715 -> 3668
$ws.Cells.Item(46, 11).Value = 12149.6  # K46: 13388.444 -> 12149.6
$ws.Cells.Item(46, 12).Value = 3668  # L46: 3396.5715 -> 3668
$ws.Cells.Item(46, 13).Value = -11961.6  # M46: -13200.444 -> -11961.6
$ws.Cells.Item(46, 14).Value = -4044  # N46: -3772.5715 -> -4044
$ws.Cells.Item(55, 8).Value = 11895.223  # H55: 8969.666999999999 -> 11895.223
$ws.Cells.Item(55, 9).Value = 929.5714  # I55: 708.6 -> 929.5714
$ws.Cells.Item(55, 11).Value = 929.5714  # K55: 708.6 -> 929.5714
$ws.Cells.Item(55, 13).Value = -756.5714  # M55: -535.6 -> -756.5714
$ws.Cells.Item(122, 8).Value = 14290593  # H122: 12504681 -> 14290593
$ws.Cells.Item(122, 9).Value = 5254  # I122: 4993.4 -> 5254
$ws.Cells.Item(122, 11).Value = 15762  # K122: 14980.2 -> 15762
$ws.Cells.Item(122, 13).Value = -13312  # M122: -12530.2 -> -13312
$ws.Cells.Item(132, 8).Value = 3345.2222  # H132: 2251.647 -> 3345.2222
$ws.Cells.Item(132, 9).Value = 2819.4  # I132: 1662.7142 -> 2819.4
$ws.Cells.Item(132, 10).Value = 4002.5  # J132: 5000 -> 4002.5
$ws.Cells.Item(132, 11).Value = 8458.200000000001  # K132: 4988.142599999999 -> 8458.200000000001
$ws.Cells.Item(132, 12).Value = 12007.5  # L132: 15000 -> 12007.5
$ws.Cells.Item(132, 13).Value = -5928.200000000001  # M132: -2458.142599999999 -> -5928.200000000001
$ws.Cells.Item(132, 14).Value = -17067.5  # N132: -20060 -> -17067.5
$ws.Cells.Item(136, 8).Value = 2378.25  # H136: 2217.5264 -> 2378.25
$ws.Cells.Item(136, 10).Value = 3241.875  # J136: 3663.3333 -> 3241.875
$ws.Cells.Item(136, 12).Value = 9725.625  # L136: 10989.9999 -> 9725.625
$ws.Cells.Item(136, 14).Value = -14825.625  # N136: -16089.9999 -> -14825.625

# ---- Sheet: WVR (33 cell updates) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(64, 8).Value = 58333.332  # H64: 55000 -> 58333.332
$ws.Cells.Item(64, 10).Value = 65000  # J64: 0 -> 65000
$ws.Cells.Item(64, 12).Value = 65000  # L64: 0 -> 65000
$ws.Cells.Item(64, 14).Value = -65496  # N64: None -> -65496
$ws.Cells.Item(67, 8).Value = 58333.332  # H67: 55000 -> 58333.332
$ws.Cells.Item(67, 10).Value = 65000  # J67: 0 -> 65000
$ws.Cells.Item(67, 12).Value = 65000  # L67: 0 -> 65000
$ws.Cells.Item(67, 14).Value = -66716  # N67: None -> -66716
$ws.Cells.Item(96, 8).Value = 4395881  # H96: 4397481 -> 4395881
$ws.Cells.Item(96, 9).Value = 12444.444  # I96: 13900 -> 12444.444
$ws.Cells.Item(96, 10).Value = 17546190  # J96: 13164642 -> 17546190
$ws.Cells.Item(96, 11).Value = 12444.444  # K96: 13900 -> 12444.444
$ws.Cells.Item(96, 12).Value = 17546190  # L96: 13164642 -> 17546190
$ws.Cells.Item(96, 13).Value = -11071.444  # M96: -12527 -> -11071.444
$ws.Cells.Item(96, 14).Value = -17548936  # N96: -13167388 -> -17548936
$ws.Cells.Item(132, 8).Value = 11375.594  # H132: 9459.154 -> 11375.594
$ws.Cells.Item(132, 9).Value = 14308.5  # I132: 11544.6 -> 14308.5
$ws.Cells.Item(132, 10).Value = 2576.875  # J132: 2507.6667 -> 2576.875
$ws.Cells.Item(132, 11).Value = 42925.5  # K132: 34633.8 -> 42925.5
$ws.Cells.Item(132, 12).Value = 7730.625  # L132: 7523.000100000001 -> 7730.625
$ws.Cells.Item(132, 13).Value = -40395.5  # M132: -32103.8 -> -40395.5
$ws.Cells.Item(132, 14).Value = -12790.625  # N132: -12583.0001 -> -12790.625
$ws.Cells.Item(136, 8).Value = 941.8182  # H136: 975.7143 -> 941.8182
$ws.Cells.Item(136, 10).Value = 0  # J136: 1100 -> 0
$ws.Cells.Item(136, 12).Value = 0  # L136: 3300 -> 0
$ws.Cells.Item(136, 14).ClearContents()  # N136 was -8400, now blank
$ws.Cells.Item(139, 8).Value = 79034.5  # H139: 79583.164 -> 79034.5
$ws.Cells.Item(139, 9).Value = 76250  # I139: 77500 -> 76250
$ws.Cells.Item(139, 10).Value = 79962.664  # J139: 79999.8 -> 79962.664
$ws.Cells.Item(139, 11).Value = 76250  # K139: 76250 -> 76250
$ws.Cells.Item(139, 12).Value = 79962.664  # L139: 79999.8 -> 79962.664
$ws.Cells.Item(139, 13).Value = -71110  # M139: -72360 -> -71110
$ws.Cells.Item(139, 14).Value = -90242.664  # N139: -90279.8 -> -90242.664
